$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (price & volume figures) scraped on
# Thu Nov 21 02:36:45 UTC 2024. Column D (Price) values are forced to
# Text format so that numeric-looking strings (e.g. "231.54") are not
# auto-converted to numbers by Excel, preserving the original formatting.

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "94.307.26"
$ws.Cells.Item(2, 5).Value = "  +2.70%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.041.69"
$ws.Cells.Item(3, 5).Value = "  -1.75%  "
$ws.Cells.Item(4, 5).Value = "  +0.06%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "231.54"
$ws.Cells.Item(5, 5).Value = "  -0.64%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "596.32"
$ws.Cells.Item(6, 5).Value = "  -2.40%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.09"
$ws.Cells.Item(7, 5).Value = "  -0.52%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.373"
$ws.Cells.Item(8, 5).Value = "  -2.84%  "
$ws.Cells.Item(9, 5).Value = "  +0.08%  "
$ws.Cells.Item(10, 2).Value = "LidoStakedEther"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "3.043.40"
$ws.Cells.Item(10, 5).Value = "  -1.52%  "
$ws.Cells.Item(11, 2).Value = "Cardano"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.772"
$ws.Cells.Item(11, 5).Value = "  +1.12%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.194"
$ws.Cells.Item(12, 5).Value = "  -1.85%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "94.158.76"
$ws.Cells.Item(13, 5).Value = "  +2.34%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.0000232"
$ws.Cells.Item(14, 5).Value = "  -4.25%  "
$ws.Cells.Item(15, 2).Value = "Toncoin"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "5.25"
$ws.Cells.Item(15, 5).Value = "  -2.64%  "
$ws.Cells.Item(16, 2).Value = "Avalanche"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "32.58"
$ws.Cells.Item(16, 5).Value = "  -2.71%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "3.616.81"
$ws.Cells.Item(17, 5).Value = "  -1.54%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "3.076.23"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "3.48"
$ws.Cells.Item(19, 5).Value = "  -7.59%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "13.94"
$ws.Cells.Item(20, 5).Value = "  -3.11%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "5.57"
$ws.Cells.Item(21, 5).Value = "  -4.20%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "433.71"
$ws.Cells.Item(22, 5).Value = "  -0.58%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "8.58"
$ws.Cells.Item(23, 5).Value = "  -5.43%  "
$ws.Cells.Item(24, 2).Value = "LEO"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "8.38"
$ws.Cells.Item(24, 5).Value = "  +3.17%  "
$ws.Cells.Item(25, 2).Value = "PEPE"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.0000182"
$ws.Cells.Item(25, 5).Value = "  -5.85%  "
$ws.Cells.Item(26, 2).Value = "NEARProtocol"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "5.39"
$ws.Cells.Item(26, 5).Value = "  -4.12%  "
$ws.Cells.Item(27, 2).Value = "Litecoin"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "83.61"
$ws.Cells.Item(27, 5).Value = "  -2.16%  "
$ws.Cells.Item(28, 2).Value = "Aptos"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "11.34"
$ws.Cells.Item(28, 5).Value = "  -0.52%  "
$ws.Cells.Item(29, 2).Value = "WrappedeETH"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "3.222.28"
$ws.Cells.Item(29, 5).Value = "  -1.43%  "
$ws.Cells.Item(30, 2).Value = "Dai"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.997"
$ws.Cells.Item(30, 5).Value = "  -0.13%  "
$ws.Cells.Item(31, 2).Value = "Hedera"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.131"
$ws.Cells.Item(31, 5).Value = "  +1.43%  "
$ws.Cells.Item(32, 2).Value = "Stellar"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.237"
$ws.Cells.Item(32, 5).Value = "  +0.96%  "
$ws.Cells.Item(33, 2).Value = "Cronos"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.171"
$ws.Cells.Item(33, 5).Value = "  -2.83%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "8.72"
$ws.Cells.Item(34, 5).Value = "  -3.15%  "
$ws.Cells.Item(35, 2).Value = "EthereumClassic"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "24.82"
$ws.Cells.Item(35, 5).Value = "  -2.71%  "
$ws.Cells.Item(36, 2).Value = "Kaspa"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.147"
$ws.Cells.Item(36, 5).Value = "  -7.13%  "
$ws.Cells.Item(37, 2).Value = "RenderToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "7.11"
$ws.Cells.Item(37, 5).Value = "  -7.76%  "
$ws.Cells.Item(38, 2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.782"
$ws.Cells.Item(38, 5).Value = "  -21.85%  "
$ws.Cells.Item(39, 2).Value = "Bittensor"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "481.46"
$ws.Cells.Item(39, 5).Value = "  +3.40%  "
$ws.Cells.Item(40, 2).Value = "WhiteBITCoin"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "24.03"
$ws.Cells.Item(40, 5).Value = "  +0.97%  "
$ws.Cells.Item(41, 2).Value = "PancakeSwap"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.82"
$ws.Cells.Item(41, 5).Value = "  -2.87%  "
$ws.Cells.Item(42, 2).Value = "MantraDAO"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.71"
$ws.Cells.Item(42, 5).Value = "  -4.91%  "
$ws.Cells.Item(43, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.421"
$ws.Cells.Item(43, 5).Value = "  -3.52%  "
$ws.Cells.Item(44, 2).Value = "Fetch.AI"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "1.19"
$ws.Cells.Item(44, 5).Value = "  -5.70%  "
$ws.Cells.Item(45, 2).Value = "USDe"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.00"
$ws.Cells.Item(45, 5).Value = "  +0.03%  "
$ws.Cells.Item(46, 2).Value = "dogwifhat"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "3.00"
$ws.Cells.Item(46, 5).Value = "  -6.87%  "
$ws.Cells.Item(47, 2).Value = "Monero"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "159.69"
$ws.Cells.Item(47, 5).Value = "  -1.78%  "
$ws.Cells.Item(48, 2).Value = "ARBITRUM"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.657"
$ws.Cells.Item(48, 5).Value = "  -3.16%  "
$ws.Cells.Item(49, 2).Value = "Stacks"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "1.76"
$ws.Cells.Item(49, 5).Value = "  -4.05%  "
$ws.Cells.Item(50, 2).Value = "FLOKI"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.000266"
$ws.Cells.Item(50, 5).Value = "  +10.44%  "
$ws.Cells.Item(51, 2).Value = "OKB"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "43.51"
$ws.Cells.Item(51, 5).Value = "  -0.80%  "
